$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"3.262296333333333"
$ws.Range("H2").Value = [double]"9.786889"
$ws.Range("I2").Value = [double]"0.01915820289899999"
$ws.Range("J2").Value = [double]"0.01915820289899999"
$ws.Range("M2").Value = [double]"35.585194"
$ws.Range("N2").Value = [double]"106.755582"
$ws.Range("O2").Value = [double]"0.9972091466993565"
$ws.Range("P2").Value = [double]"0.9972091466993567"
$ws.Range("Q2").Value = [double]"116.0894479071553"
$ws.Range("R2").Value = [double]"1044.805031164398"
$ws.Range("S2").Value = [double]"0.01910473516520492"
$ws.Range("T2").Value = [double]"0.01910473516520492"
$ws.Range("G3").Value = [double]"3.262296333333333"
$ws.Range("H3").Value = [double]"9.786889"
$ws.Range("I3").Value = [double]"0.01915820289899999"
$ws.Range("J3").Value = [double]"0.01915820289899999"
$ws.Range("M3").Value = [double]"0.093901"
$ws.Range("N3").Value = [double]"0.281703"
$ws.Range("O3").Value = [double]"0.002631401590341653"
$ws.Range("P3").Value = [double]"0.002631401590341654"
$ws.Range("Q3").Value = [double]"0.3063328879963333"
$ws.Range("R3").Value = [double]"2.756995991967"
$ws.Range("S3").Value = [double]"5.041292557651665E-05"
$ws.Range("T3").Value = [double]"5.041292557651666E-05"
$ws.Range("G4").Value = [double]"3.262296333333333"
$ws.Range("H4").Value = [double]"9.786889"
$ws.Range("I4").Value = [double]"0.01915820289899999"
$ws.Range("J4").Value = [double]"0.01915820289899999"
$ws.Range("K4").Value = [double]"1"
$ws.Range("L4").Value = [double]"0.3333333333333333"
$ws.Range("M4").Value = [double]"0.00569"
$ws.Range("N4").Value = [double]"0.01707"
$ws.Range("O4").Value = [double]"0.0001594517103017434"
$ws.Range("P4").Value = [double]"0.0001594517103017434"
$ws.Range("Q4").Value = [double]"0.01856246613666667"
$ws.Range("R4").Value = [double]"0.16706219523"
$ws.Range("S4").Value = [double]"3.054808218553367E-06"
$ws.Range("T4").Value = [double]"3.054808218553368E-06"
$ws.Range("I5").Value = [double]"0.8527862647199704"
$ws.Range("J5").Value = [double]"0.8527862647199704"
$ws.Range("M5").Value = [double]"35.585194"
$ws.Range("N5").Value = [double]"106.755582"
$ws.Range("O5").Value = [double]"0.9972091466993565"
$ws.Range("P5").Value = [double]"0.9972091466993567"
$ws.Range("Q5").Value = [double]"5167.472501260236"
$ws.Range("R5").Value = [double]"46507.25251134213"
$ws.Range("S5").Value = [double]"0.8504062633583332"
$ws.Range("T5").Value = [double]"0.8504062633583334"
$ws.Range("I6").Value = [double]"0.8527862647199704"
$ws.Range("J6").Value = [double]"0.8527862647199704"
$ws.Range("M6").Value = [double]"0.093901"
$ws.Range("N6").Value = [double]"0.281703"
$ws.Range("O6").Value = [double]"0.002631401590341653"
$ws.Range("P6").Value = [double]"0.002631401590341654"
$ws.Range("Q6").Value = [double]"13.63575073781633"
$ws.Range("R6").Value = [double]"122.721756640347"
$ws.Range("S6").Value = [double]"0.002244023133205648"
$ws.Range("T6").Value = [double]"0.002244023133205649"
$ws.Range("I7").Value = [double]"0.8527862647199704"
$ws.Range("J7").Value = [double]"0.8527862647199704"
$ws.Range("K7").Value = [double]"1"
$ws.Range("L7").Value = [double]"0.3333333333333333"
$ws.Range("M7").Value = [double]"0.00569"
$ws.Range("N7").Value = [double]"0.01707"
$ws.Range("O7").Value = [double]"0.0001594517103017434"
$ws.Range("P7").Value = [double]"0.0001594517103017434"
$ws.Range("Q7").Value = [double]"0.8262683219366668"
$ws.Range("R7").Value = [double]"7.436414897430001"
$ws.Range("S7").Value = [double]"0.0001359782284314346"
$ws.Range("T7").Value = [double]"0.0001359782284314346"
$ws.Range("G8").Value = [double]"21.305189"
$ws.Range("H8").Value = [double]"63.915567"
$ws.Range("I8").Value = [double]"0.1251171236325075"
$ws.Range("J8").Value = [double]"0.1251171236325075"
$ws.Range("M8").Value = [double]"35.585194"
$ws.Range("N8").Value = [double]"106.755582"
$ws.Range("O8").Value = [double]"0.9972091466993565"
$ws.Range("P8").Value = [double]"0.9972091466993567"
$ws.Range("Q8").Value = [double]"758.149283771666"
$ws.Range("R8").Value = [double]"6823.343553944994"
$ws.Range("S8").Value = [double]"0.1247679400950507"
$ws.Range("T8").Value = [double]"0.1247679400950508"
$ws.Range("G9").Value = [double]"21.305189"
$ws.Range("H9").Value = [double]"63.915567"
$ws.Range("I9").Value = [double]"0.1251171236325075"
$ws.Range("J9").Value = [double]"0.1251171236325075"
$ws.Range("M9").Value = [double]"0.093901"
$ws.Range("N9").Value = [double]"0.281703"
$ws.Range("O9").Value = [double]"0.002631401590341653"
$ws.Range("P9").Value = [double]"0.002631401590341654"
$ws.Range("Q9").Value = [double]"2.000578552289"
$ws.Range("R9").Value = [double]"18.005206970601"
$ws.Range("S9").Value = [double]"0.0003292333981055536"
$ws.Range("T9").Value = [double]"0.0003292333981055537"
$ws.Range("G10").Value = [double]"21.305189"
$ws.Range("H10").Value = [double]"63.915567"
$ws.Range("I10").Value = [double]"0.1251171236325075"
$ws.Range("J10").Value = [double]"0.1251171236325075"
$ws.Range("K10").Value = [double]"1"
$ws.Range("L10").Value = [double]"0.3333333333333333"
$ws.Range("M10").Value = [double]"0.00569"
$ws.Range("N10").Value = [double]"0.01707"
$ws.Range("O10").Value = [double]"0.0001594517103017434"
$ws.Range("P10").Value = [double]"0.0001594517103017434"
$ws.Range("Q10").Value = [double]"0.12122652541"
$ws.Range("R10").Value = [double]"1.09103872869"
$ws.Range("S10").Value = [double]"1.995013935123801E-05"
$ws.Range("T10").Value = [double]"1.995013935123801E-05"
$ws.Range("G11").Value = [double]"0.500358"
$ws.Range("H11").Value = [double]"1.501074"
$ws.Range("I11").Value = [double]"0.002938408748521978"
$ws.Range("J11").Value = [double]"0.002938408748521978"
$ws.Range("M11").Value = [double]"35.585194"
$ws.Range("N11").Value = [double]"106.755582"
$ws.Range("O11").Value = [double]"0.9972091466993565"
$ws.Range("P11").Value = [double]"0.9972091466993567"
$ws.Range("Q11").Value = [double]"17.805336499452"
$ws.Range("R11").Value = [double]"160.248028495068"
$ws.Range("S11").Value = [double]"0.002930208080767526"
$ws.Range("T11").Value = [double]"0.002930208080767526"
$ws.Range("G12").Value = [double]"0.500358"
$ws.Range("H12").Value = [double]"1.501074"
$ws.Range("I12").Value = [double]"0.002938408748521978"
$ws.Range("J12").Value = [double]"0.002938408748521978"
$ws.Range("M12").Value = [double]"0.093901"
$ws.Range("N12").Value = [double]"0.281703"
$ws.Range("O12").Value = [double]"0.002631401590341653"
$ws.Range("P12").Value = [double]"0.002631401590341654"
$ws.Range("Q12").Value = [double]"0.046984116558"
$ws.Range("R12").Value = [double]"0.422857049022"
$ws.Range("S12").Value = [double]"7.73213345393456E-06"
$ws.Range("T12").Value = [double]"7.732133453934562E-06"
$ws.Range("G13").Value = [double]"0.500358"
$ws.Range("H13").Value = [double]"1.501074"
$ws.Range("I13").Value = [double]"0.002938408748521978"
$ws.Range("J13").Value = [double]"0.002938408748521978"
$ws.Range("K13").Value = [double]"1"
$ws.Range("L13").Value = [double]"0.3333333333333333"
$ws.Range("M13").Value = [double]"0.00569"
$ws.Range("N13").Value = [double]"0.01707"
$ws.Range("O13").Value = [double]"0.0001594517103017434"
$ws.Range("P13").Value = [double]"0.0001594517103017434"
$ws.Range("Q13").Value = [double]"0.00284703702"
$ws.Range("R13").Value = [double]"0.02562333318"
$ws.Range("S13").Value = [double]"4.685343005174348E-07"
$ws.Range("T13").Value = [double]"4.685343005174349E-07"
